$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume change (E) columns
$ws.Range("D2").Value = "31.042.36"
$ws.Range("E2").Value = "  +3.81%  "
$ws.Range("D3").Value = "1.691.84"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.75"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.533"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.46"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "1.937.42"
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("D13").Value = "1.680.01"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.29"
$ws.Range("E14").Value = "  +8.22%  "
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.13"
$ws.Range("E16").Value = "  +6.30%  "
$ws.Range("D17").Value = "31.051.95"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.90"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.64"
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.09"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.25"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.94"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.73"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.57"
$ws.Range("E31").Value = "  +5.34%  "
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("E33").Value = "  +4.80%  "
$ws.Range("D34").Value = "1.516.94"
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "83.12"
$ws.Range("E37").Value = "  +9.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.615"
$ws.Range("E38").Value = "  +9.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0179"
$ws.Range("E39").Value = "  +4.24%  "
$ws.Range("E40").Value = "  -4.17%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.846"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.04"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "51.88"
$ws.Range("E47").Value = "  +6.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.57"
$ws.Range("E48").Value = "  +4.37%  "
$ws.Range("D49").Value = "1.824.24"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  +9.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "93.75"
$ws.Range("E51").Value = "  +1.07%  "
